$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on all touched cells so numeric-looking values
# (e.g. "226.97") stay stored as text, matching the source data feed
# (inline strings), not auto-converted to numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "39.776.80"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.45%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.161.00"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.66%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.97"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.626"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.57%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "63.30"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.73%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0843"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.20%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.03%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.87"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.482.60"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.83"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.805"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.47%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.50"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.19%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.157.82"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +4.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "39.684.69"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.73"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.00"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.79%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "229.68"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.93%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.46%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.31"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -7.45%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "172.47"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.52"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.38%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.03%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.35%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.79"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.22%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +6.38%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.13%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.68"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.71%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.62%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.37%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.70"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +5.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.40"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.70%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +22.21%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "102.44"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.93%  "
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0227"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.75%  "
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.71"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.79%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.514.91"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.20"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.46%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.78%  "
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0921"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.83%  "
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "HuobiToken"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.80"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.04%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "49.91"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +8.44%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.99"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.93%  "
